$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# Remove the "Smiley Face 1" shape entirely.
$smiley = $s.Shapes.Item("Smiley Face 1")
$smiley.Delete()

# Reposition / resize the "TextBox 2" shape (EMU -> points, 12700 EMU = 1 pt).
# A tiny half-EMU epsilon counters the host's truncation (rather than
# rounding) when it converts the Single-precision point value back to EMU.
$textBox = $s.Shapes.Item("TextBox 2")
$emu = 12700.0
$eps = 0.5 / $emu
$textBox.Left   = (1740024 / $emu) + $eps
$textBox.Top    = (2787588 / $emu) + $eps
$textBox.Width  = (8442664 / $emu) + $eps
$textBox.Height = (923330 / $emu) + $eps

# Bump the "THANK YOU!!!!!" run's font size from 32pt to 54pt.
$tr = $textBox.TextFrame.TextRange
$fullText = $tr.Text
$idx = $fullText.IndexOf("THANK YOU!!!!!")
$run = $tr.Characters($idx + 1, "THANK YOU!!!!!".Length)
$run.Font.Size = 54
